$wb = $excel.ActiveWorkbook

# xlPasteValues — used to "flatten" a formula-computed text result into a
# plain stored value without Excel's autodetection re-parsing an ISO-looking
# date string (e.g. "2026-02-17") back into a date serial number.
$xlPasteValues = -4163

# --- Sheet "Summary": update aggregate metrics (trade #10 just closed) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.87   # Current Capital
$summary.Range("B4").Value = -0.13     # Total P&L $
$summary.Range("B5").Value = -0.26     # Total P&L %
$summary.Range("B6").Value = 10        # Total Trades
$summary.Range("B8").Value = 5         # Losing Trades
$summary.Range("B9").Value = 30        # Win Rate %

# --- Sheet "Strategy Status": update MarketMaking row (row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.87
$status.Range("D4").Value = 10
$status.Range("E4").Value = -0.13
$status.Range("F4").Value = -0.13
$status.Range("G4").Value = 30

# --- Append new trade #10 row to both "All Trades" and "MarketMaking" ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($name in $tradeSheets) {
    $ws = $wb.Worksheets.Item($name)
    $row = 11

    $ws.Cells.Item($row, 1).Value = 10

    # Date column: route through a text formula + paste-values so the literal
    # string "2026-02-17" is stored as text (matching the other rows) instead
    # of being auto-converted to a date serial number.
    $ws.Cells.Item($row, 2).Formula = '="2026-02-17"'
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($row, 3).Value = "15:14:22"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.03
    $ws.Cells.Item($row, 7).Value = 0.02
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -33.3333
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 99.87
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.15
}

$excel.CutCopyMode = $false
